$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.812.67"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").Value = "1.541.67"
$ws.Range("E3").Value = "  -0.85%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "289.78"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3909"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3206"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.54%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.41"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07173"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.063"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.65%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.614"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.54"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -7.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.625"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "1.543.66"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001098"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06561"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.27"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.126"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.28"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.96"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.37%  "
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").Value = "21.823.23"
$ws.Range("E25").Value = "  -1.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.370"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -6.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "145.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.46"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -3.43%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.841"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.47%  "
$ws.Range("D30").Value = "1.714.70"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "117.31"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -3.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9752"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -10.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.855"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08186"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.933"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.583"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -14.89%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06034"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02235"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -4.25%  "
$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.087"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2031"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -6.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.184"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.005"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5751"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.751"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.94"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5533"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "116.49"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.863"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.128"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06758"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.67%  "
